# Bitacora Laboratorio - Parcial 04 - actualizacion 14 abril 2024
# 1) Renombra la hoja original a "Concentrado" y agrega las columnas de
#    evaluacion P7-P9 con sus puntajes.
# 2) Duplica la hoja original para crear "Asistencia", con fechas de
#    sesion y marcas de asistencia.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Crea la hoja "Asistencia" como copia de la hoja original, ANTES de
#     modificar la hoja original, para heredar alumnos/formatos tal cual. ---
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Asistencia"

# Renombra la hoja original
$ws1.Name = "Concentrado"

# ======================================================================
# Hoja "Concentrado": nuevas columnas de evaluacion (P7 a P9)
# ======================================================================

$headers = @{
  "E1" = "P7_Reporte"
  "F1" = "P8_Encuadre"
  "G1" = "P8_Montaje"
  "H1" = "P8_Discusión"
  "I1" = "P8_Reporte"
  "J1" = "P9 Encuadre"
  "K1" = "P9 Montaje"
  "L1" = "P9_Discusión"
  "M1" = "P9_Reporte"
  "N1" = "Puntaje"
  "O1" = "Calificación"
}
foreach ($addr in $headers.Keys) {
  $cell = $ws1.Range($addr)
  $cell.Value = $headers[$addr]
  $cell.Font.Bold = $true
  $cell.HorizontalAlignment = -4108
}

# Puntajes por alumno (filas 2 a 15), columnas F a L
$scores = @{
  2  = @{ F=5; G=5; H=5; I=9; J=5; K=5 }
  3  = @{ F=5; G=5; H=5; I=7; J=5; K=5; L=5 }
  4  = @{ F=5; G=5; H=5; I=9; J=5; K=5 }
  5  = @{ F=5; G=5; H=5; I=7; J=0; K=0; L=0 }
  6  = @{ F=5; G=5; H=5;         J=0; K=5; L=5 }
  7  = @{ F=5; G=5; H=5;         J=5; K=5; L=5 }
  8  = @{ F=5; G=5; H=5;         J=5; K=5; L=5 }
  9  = @{ F=5; G=5; H=5; I=9; J=5; K=5 }
  10 = @{ F=4; G=5; H=5;         J=5; K=5 }
  11 = @{ F=5; G=5; H=5;         J=5; K=5; L=5 }
  12 = @{ F=5; G=5; H=5; I=7; J=5; K=5; L=5 }
  13 = @{ F=5; G=5; H=5; I=9; J=5; K=5 }
  14 = @{ F=5; G=5; H=5;         J=5; K=5 }
  15 = @{ F=5; G=5; H=5;         J=5; K=5 }
}

foreach ($row in $scores.Keys) {
  $cols = $scores[$row]
  foreach ($col in $cols.Keys) {
    $addr = "$col$row"
    $cell = $ws1.Range($addr)
    $cell.Value = $cols[$col]
    $cell.HorizontalAlignment = -4108
  }
}

$ws1.Columns.Item(5).EntireColumn.AutoFit()
$ws1.Columns.Item(6).EntireColumn.AutoFit()
$ws1.Columns.Item(7).EntireColumn.AutoFit()
$ws1.Columns.Item(8).EntireColumn.AutoFit()
$ws1.Columns.Item(9).EntireColumn.AutoFit()
$ws1.Columns.Item(10).EntireColumn.AutoFit()
$ws1.Columns.Item(11).EntireColumn.AutoFit()
$ws1.Columns.Item(12).EntireColumn.AutoFit()
$ws1.Columns.Item(13).EntireColumn.AutoFit()
$ws1.Columns.Item(14).EntireColumn.AutoFit()
$ws1.Columns.Item(15).EntireColumn.AutoFit()

$ws1.Range("A1:D15").Select()

# ======================================================================
# Hoja "Asistencia": fechas de sesion y marcas de asistencia
# ======================================================================

# Limpia los encabezados heredados (P8_Avance / P8_Discusión) de E1:F1
$ws2.Range("E1:F1").Clear()

$dates = @{
  "E1" = 45357
  "F1" = 45364
  "G1" = 45371
  "H1" = 45385
  "I1" = 45392
}
foreach ($addr in $dates.Keys) {
  $cell = $ws2.Range($addr)
  $cell.Value = $dates[$addr]
  $cell.NumberFormat = "d-mmm"
  $cell.HorizontalAlignment = -4108
}

$ws2.Range("H5").Value = 1
$ws2.Range("H5").HorizontalAlignment = -4108
$ws2.Range("I5").Value = 1
$ws2.Range("I5").HorizontalAlignment = -4108

$ws2.Columns.Item(5).EntireColumn.AutoFit()
$ws2.Columns.Item(6).EntireColumn.AutoFit()
$ws2.Columns.Item(7).EntireColumn.AutoFit()
$ws2.Columns.Item(8).EntireColumn.AutoFit()
$ws2.Columns.Item(9).EntireColumn.AutoFit()

$ws2.Range("G5").Select()
$ws2.Activate()

Write-Host "Edit complete"
